$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Data table layout: col A = Date, col B = Non-HTTPS URLs, col C = HTTPS
# URLs, one header row followed by one row per day. This export rolls the
# whole window forward by one day: the oldest day (row 2) is dropped, every
# remaining day's counts shift up a row, and a brand-new day (no data yet)
# is appended after the previous last row.
$usedRows = $ws.UsedRange.Rows.Count()
$lastRow = $usedRows

$lastDateStr = $ws.Cells.Item($lastRow, 1).Value()
$lastDate = [datetime]::ParseExact($lastDateStr, "yyyy-MM-dd", $null)
$newDateStr = $lastDate.AddDays(1).ToString("yyyy-MM-dd")

# Iterate ascending from the top: row r is rewritten using row r+1's
# CURRENT values, which are still untouched at that point since r+1 > r
# hasn't been visited by this loop yet.
for ($r = 2; $r -le ($lastRow - 1); $r++) {
    $srcDate = $ws.Cells.Item($r + 1, 1).Value()
    $srcB = $ws.Cells.Item($r + 1, 2).Value()
    $srcC = $ws.Cells.Item($r + 1, 3).Value()

    # A plain .Value assignment of a date-shaped string ("2025-09-01") gets
    # silently reinterpreted as a date serial number by the host's type
    # inference, which would flip the cell away from the plain
    # shared-string text the rest of the column uses (and pull in a new
    # number-format style). Forcing text entry with a leading apostrophe
    # keeps it a string, then ClearFormats drops the quote-prefix flag that
    # leaves behind so the cell's style matches its neighbours exactly.
    $ws.Cells.Item($r, 1).Value = "'" + $srcDate
    $ws.Cells.Item($r, 1).ClearFormats()

    $ws.Cells.Item($r, 2).Value = $srcB
    $ws.Cells.Item($r, 3).Value = $srcC
}

$ws.Cells.Item($lastRow, 1).Value = "'" + $newDateStr
$ws.Cells.Item($lastRow, 1).ClearFormats()
$ws.Cells.Item($lastRow, 2).Value = 0
$ws.Cells.Item($lastRow, 3).Value = 0
